$d = $word.ActiveDocument

# Insert a new, centered contact-info paragraph directly after the name
# ("Dheeraj Chand") and before the "PROFESSIONAL SUMMARY" heading, by
# replacing "Dheeraj Chand" with itself plus a paragraph break and the
# contact line. This keeps the new run free of the name run's direct
# bold/size formatting (which InsertParagraphAfter would otherwise
# inherit) while still picking up the paragraph's centered alignment.
$d.Content.Find.Execute(
    "Dheeraj Chand", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX",
    2)
